# Documentação do Projeto / especificacao-funcional-estruturada.xlsx
# -------------------------------------------------------------------
# 1) Fix the typo "Alertas via Slak" -> "Alertas via Slack" (B6), keeping
#    the original rich-text layout: a bold run for the title followed by a
#    plain run holding the trailing line break.
# 2) Replace the stale "Cadastro e gerenciamento de métricas..." requirement
#    text in C7 with the new "captura de dados do computador..." text.
# 3) Re-fit the row heights that reflect the new wrapped text, and restore
#    the selection that was left on the sheet.

# Shared-string table entries are (re)created in first-use order when the
# workbook is re-serialised, so the C7 replacement is written before the B6
# rich-text replacement to land the two new strings at the same shared-string
# slots the original edit produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) C7: old "Cadastro e gerenciamento..." -> new "captura de dados..." -
$ws.Range("C7").Value = "captura de dados do computador de forma esporádica, conforme período determinado pela empresa."

# --- 2) B6: "Alertas via Slak" -> "Alertas via Slack" ---------------------
$ws.Range("B6").Value = "Alertas via Slack" + [char]10

$titleRun = $ws.Range("B6").Characters(1, 17)
$titleRun.Font.Bold = $true
$titleRun.Font.Name = "Arial"
$titleRun.Font.Size = 12
$titleRun.Font.ColorIndex = -4105

$breakRun = $ws.Range("B6").Characters(18, 1)
$breakRun.Font.Bold = $false
$breakRun.Font.Name = "Arial"
$breakRun.Font.Size = 12
$breakRun.Font.ColorIndex = -4105

# --- 3) Row heights reflecting the edited wrapped text ---------------------
$ws.Rows.Item(6).RowHeight = 165.6
$ws.Rows.Item(9).RowHeight = 151.8
$ws.Rows.Item(11).RowHeight = 69

# --- 4) Selection left on the sheet after editing --------------------------
$ws.Range("K7").Select()

Write-Output "edits applied"
